# Apply the changes described by the commit "Changes in Shiny Application":
#  1. Fix the shared string "I am gloomy about my<TAB>future" -> "I am gloomy about my future"
#     (a stray tab character between "my" and "future" is replaced with a space).
#     Every cell in column J that used this text has to be rewritten with the
#     corrected text so the shared string table collapses back down to a single
#     (fixed) entry instead of growing a duplicate.
#  2. Update the sheet view: the user scrolled down (topLeftCell around row 10)
#     and the final selection left on the sheet was J51.
#  3. Column J (10) was widened to fit the long question text (customWidth).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Correct the mis-typed shared string (tab -> space) everywhere it is used ---
$fixedText = "I am gloomy about my future"
$cellsToFix = @("J9","J14","J16","J18","J20","J24","J34","J39","J40","J42","J45","J46","J47","J50")
foreach ($addr in $cellsToFix) {
    $ws.Range($addr).Value = $fixedText
}

# --- 3. Widen column J (10) to fit the question text ---
$ws.Columns.Item(10).ColumnWidth = 75.65

# --- 2. Scroll the view down and leave the final selection on J51 ---
$ws.Activate() | Out-Null
$ws.Range("A10").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J51").Select() | Out-Null
